$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md) now ready for handoff ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-03 04:41:43"

# --- zh-cn sheet: row 3 (b.md) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text entry (avoids auto-coercion to a Boolean),
# then resetting the style clears the resulting quote-prefix flag.
$ws.Range("F3").Value = "'False"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-03 04:41:35"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce15f614f1085b4ba376681c5a9d6604108d870b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a257fd6e8c5cca670fa53106c5f5e0ec2d6a4c18/e2e/b.md."
# Column P (Error Detail) widens to match the other long-text columns (G/J = 40).
$ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# --- de-de sheet: row 3 (b.md) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'False"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-09-03 04:41:43"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce15f614f1085b4ba376681c5a9d6604108d870b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a257fd6e8c5cca670fa53106c5f5e0ec2d6a4c18/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(7).ColumnWidth
